# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 146 ("Berenjena" sheet), shifting
# the existing rows 146-188 down to 147-189 (dimension grows from R188 to R189).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 146, pushing all following rows down.
$ws.Rows.Item(146).Insert()

# Populate the newly inserted row 146 with the new weekly record.
$ws.Cells.Item(146, 1).Value  = 8
$ws.Cells.Item(146, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(146, 3).Value  = "Coquimbo"
$ws.Cells.Item(146, 4).Value  = 44932
$ws.Cells.Item(146, 5).Value  = 4
$ws.Cells.Item(146, 6).Value  = 100112001
$ws.Cells.Item(146, 7).Value  = "Berenjena"
$ws.Cells.Item(146, 8).Value  = "Sin especificar"
$ws.Cells.Item(146, 9).Value  = "Primera"
$ws.Cells.Item(146, 10).Value = 400
$ws.Cells.Item(146, 11).Value = 11000
$ws.Cells.Item(146, 12).Value = 12000
$ws.Cells.Item(146, 13).Value = 11500
$ws.Cells.Item(146, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(146, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(146, 16).Value = 288
$ws.Cells.Item(146, 17).Value = 40
$ws.Cells.Item(146, 18).Value = "Hortaliza"
